# Re-orders the species-occurrence records in rows 52-67 of the "Artfynd"
# worksheet so that column A (Id) values match the canonical source order,
# rounds the Ost/Nord (Q/R) coordinates to whole metres, clears the now-unused
# Starttid/Sluttid (Z/AB) time cells for rows 52-66, and keeps the species-specific
# Ålder-Stadium/Kön/Aktivitet/Metod (K/L/M/N) + Publik kommentar (AC) cells in sync
# with whichever record now occupies each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targets = @(
    @{ Row=52; A=111901548; B=56398; E=100109; F="Tretåig hackspett"; G="Picoides tridactylus"; H="(Linnaeus, 1758)"; Q=477476; R=7033385; KLMN=$true; AC="ringhack äldre" }
    @{ Row=53; A=111901519; B=86223; E=4412; F="Äggvaxskivling"; G="Hygrophorus karstenii"; H="Sacc. & Cub."; Q=477765; R=7033404; KLMN=$false; AC=$null }
    @{ Row=54; A=111901518; B=86223; E=4412; F="Äggvaxskivling"; G="Hygrophorus karstenii"; H="Sacc. & Cub."; Q=477674; R=7033500; KLMN=$false; AC=$null }
    @{ Row=55; A=111901550; B=56398; E=100109; F="Tretåig hackspett"; G="Picoides tridactylus"; H="(Linnaeus, 1758)"; Q=477473; R=7033404; KLMN=$true; AC="ringhack äldre" }
    @{ Row=56; A=111901549; B=56398; E=100109; F="Tretåig hackspett"; G="Picoides tridactylus"; H="(Linnaeus, 1758)"; Q=477464; R=7033364; KLMN=$true; AC="ringhack färska" }
    @{ Row=57; A=111901585; B=56398; E=100109; F="Tretåig hackspett"; G="Picoides tridactylus"; H="(Linnaeus, 1758)"; Q=478339; R=7035076; KLMN=$true; AC="ringhack äldre" }
    @{ Row=58; A=111901551; B=56398; E=100109; F="Tretåig hackspett"; G="Picoides tridactylus"; H="(Linnaeus, 1758)"; Q=477433; R=7033429; KLMN=$true; AC="ringhack" }
    @{ Row=59; A=111901544; B=56398; E=100109; F="Tretåig hackspett"; G="Picoides tridactylus"; H="(Linnaeus, 1758)"; Q=477639; R=7033515; KLMN=$true; AC="ringhack äldre" }
    @{ Row=60; A=111901618; B=85062; E=249278; F="Barrviolspindling"; G="Cortinarius harcynicus"; H="(Pers.) M.M.Moser"; Q=477471; R=7033412; KLMN=$false; AC=$null }
    @{ Row=61; A=111901545; B=56398; E=100109; F="Tretåig hackspett"; G="Picoides tridactylus"; H="(Linnaeus, 1758)"; Q=477667; R=7033500; KLMN=$true; AC="ringhack äldre" }
    @{ Row=62; A=111901547; B=56398; E=100109; F="Tretåig hackspett"; G="Picoides tridactylus"; H="(Linnaeus, 1758)"; Q=477524; R=7033330; KLMN=$true; AC="ringhack" }
    @{ Row=63; A=111901584; B=56398; E=100109; F="Tretåig hackspett"; G="Picoides tridactylus"; H="(Linnaeus, 1758)"; Q=478211; R=7035067; KLMN=$true; AC="ringhack" }
    @{ Row=64; A=111901546; B=56398; E=100109; F="Tretåig hackspett"; G="Picoides tridactylus"; H="(Linnaeus, 1758)"; Q=477668; R=7033374; KLMN=$true; AC="ringhack äldre" }
    @{ Row=65; A=111901619; B=85062; E=249278; F="Barrviolspindling"; G="Cortinarius harcynicus"; H="(Pers.) M.M.Moser"; Q=478523; R=7034651; KLMN=$false; AC=$null }
    @{ Row=66; A=111901587; B=56543; E=103021; F="Talltita"; G="Poecile montanus"; H="(Conrad von Baldenstein, 1827)"; Q=477611; R=7033311; KLMN=$true; AC=$null }
    @{ Row=67; A=112102606; B=88899; E=3286; F="Flattoppad klubbsvamp"; G="Clavariadelphus truncatus"; H="(Quél.) Donk"; Q=478088; R=7035319; KLMN=$true; AC=$null }
)

foreach ($t in $targets) {
    $row = $t.Row

    $ws.Range("A$row").Value = $t.A
    $ws.Range("B$row").Value = $t.B
    $ws.Range("E$row").Value = $t.E

    $ws.Range("F$row").NumberFormat = "@"
    $ws.Range("F$row").Value = $t.F
    $ws.Range("G$row").NumberFormat = "@"
    $ws.Range("G$row").Value = $t.G
    $ws.Range("H$row").NumberFormat = "@"
    $ws.Range("H$row").Value = $t.H

    $ws.Range("Q$row").Value = $t.Q
    $ws.Range("R$row").Value = $t.R

    if ($t.KLMN) {
        $ws.Range("K$row").NumberFormat = "@"
        $ws.Range("K$row").Value = ""
        $ws.Range("L$row").NumberFormat = "@"
        $ws.Range("L$row").Value = ""
        $ws.Range("M$row").NumberFormat = "@"
        $ws.Range("M$row").Value = ""
        $ws.Range("N$row").NumberFormat = "@"
        $ws.Range("N$row").Value = ""
    } else {
        $ws.Range("K$row").ClearContents()
        $ws.Range("L$row").ClearContents()
        $ws.Range("M$row").ClearContents()
        $ws.Range("N$row").ClearContents()
    }

    if ($t.AC -ne $null) {
        $ws.Range("AC$row").NumberFormat = "@"
        $ws.Range("AC$row").Value = $t.AC
    } else {
        $ws.Range("AC$row").ClearContents()
    }

    # Starttid/Sluttid only ever held the literal "00:00" placeholder for
    # these rows; the source export dropped both columns for rows 52-66.
    if ($row -le 66) {
        $ws.Range("Z$row").ClearContents()
        $ws.Range("AB$row").ClearContents()
    }
}
